$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '37.832.01'
$ws.Range("E2").Value = '  +0.14%  '
$ws.Range("D3").Value = '2.089.24'
$ws.Range("E3").Value = '  +0.34%  '
$ws.Range("E4").Value = '  +0.04%  '
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = '233.82'
$c.NumberFormat = "General"
$ws.Range("E5").Value = '  -0.35%  '
$ws.Range("E6").Value = '  -0.21%  '
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = '58.44'
$c.NumberFormat = "General"
$ws.Range("E7").Value = '  -0.57%  '
$ws.Range("E9").Value = '  +0.71%  '
$ws.Range("E10").Value = '  -0.75%  '
$ws.Range("E11").Value = '  +2.81%  '
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = '15.29'
$c.NumberFormat = "General"
$ws.Range("E12").Value = '  +3.40%  '
$ws.Range("D13").Value = '2.396.70'
$ws.Range("E13").Value = '  +0.25%  '
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = '21.25'
$c.NumberFormat = "General"
$ws.Range("E14").Value = '  +0.60%  '
$ws.Range("E15").Value = '  +0.94%  '
$ws.Range("E16").Value = '  +1.13%  '
$ws.Range("D17").Value = '2.097.69'
$ws.Range("E17").Value = '  +1.09%  '
$ws.Range("D18").Value = '37.825.26'
$ws.Range("E18").Value = '  +0.32%  '
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = '6.15'
$c.NumberFormat = "General"
$ws.Range("E19").Value = '  -0.16%  '
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = '71.07'
$c.NumberFormat = "General"
$ws.Range("E20").Value = '  -0.44%  '
$ws.Range("D21").Value = '0.0₃0837'
$ws.Range("E21").Value = '  +0.20%  '
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = '229.88'
$c.NumberFormat = "General"
$ws.Range("E22").Value = '  +0.43%  '
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = '0.999'
$c.NumberFormat = "General"
$ws.Range("E23").Value = '  -0.03%  '
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = '2.41'
$c.NumberFormat = "General"
$ws.Range("E24").Value = '  +0.03%  '
$ws.Range("E25").Value = '  -1.26%  '
$ws.Range("E26").Value = '  +8.60%  '
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = '171.67'
$c.NumberFormat = "General"
$ws.Range("E27").Value = '  +1.43%  '
$ws.Range("E28").Value = '  -3.16%  '
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = '19.53'
$c.NumberFormat = "General"
$ws.Range("E29").Value = '  -0.13%  '
$ws.Range("E30").Value = '  -0.29%  '
$ws.Range("E31").Value = '  -0.03%  '
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = '4.70'
$c.NumberFormat = "General"
$ws.Range("E32").Value = '  +0.11%  '
$ws.Range("E33").Value = '  -0.09%  '
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = '4.64'
$c.NumberFormat = "General"
$ws.Range("E34").Value = '  -0.72%  '
$ws.Range("E35").Value = '  +0.32%  '
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = '1.82'
$c.NumberFormat = "General"
$ws.Range("E36").Value = '  -0.46%  '
$ws.Range("E37").Value = '  -1.92%  '
$ws.Range("E38").Value = '  -0.16%  '
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = '5.40'
$c.NumberFormat = "General"
$ws.Range("E39").Value = '  -0.64%  '
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = '0.0235'
$c.NumberFormat = "General"
$ws.Range("E40").Value = '  +8.87%  '
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = '101.47'
$c.NumberFormat = "General"
$ws.Range("E41").Value = '  +2.87%  '
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = '0.0973'
$c.NumberFormat = "General"
$ws.Range("E42").Value = '  -0.67%  '
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = '1.20'
$c.NumberFormat = "General"
$ws.Range("E44").Value = '  +2.66%  '
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = '16.79'
$c.NumberFormat = "General"
$ws.Range("E45").Value = '  +2.39%  '
$ws.Range("D46").Value = '1.454.52'
$ws.Range("E46").Value = '  -0.57%  '
$ws.Range("E47").Value = '  -0.32%  '
$ws.Range("E48").Value = '  -4.82%  '
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = '7.22'
$c.NumberFormat = "General"
$ws.Range("E49").Value = '  -3.74%  '
$ws.Range("E50").Value = '  -1.68%  '
$ws.Range("D51").Value = '2.281.10'
$ws.Range("E51").Value = '  +0.25%  '
